# Apply "did again some tests" edit to conditions_temp_test.xlsx
# Sets B:F values to 0 for several row ranges, and to doubled magnitude
# (-20,-20,40,-10,10) for another row range, on the "Hydraulic" sheet.
# Also updates the active cell/selection and scroll position of that sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hydraulic")

# Row ranges whose B:F values become 0
$zeroRanges = @(
    @(409, 418),
    @(431, 440),
    @(455, 464),
    @(480, 489)
)

foreach ($r in $zeroRanges) {
    $startRow = $r[0]
    $endRow = $r[1]
    $rng = $ws.Range("B$startRow`:F$endRow")
    $rng.Value = 0
}

# Row range 465-479 gets doubled-magnitude values: -20, -20, 40, -10, 10
for ($row = 465; $row -le 479; $row++) {
    $ws.Cells.Item($row, 2).Value = -20
    $ws.Cells.Item($row, 3).Value = -20
    $ws.Cells.Item($row, 4).Value = 40
    $ws.Cells.Item($row, 5).Value = -10
    $ws.Cells.Item($row, 6).Value = 10
}

# Update the sheet view: scroll position and selection
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 478
$ws.Range("B465:F479").Select()
$excel.ActiveCell = $ws.Range("F465")
